$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 4 ----
$ws.Range("A4").Value = 112141528
$ws.Range("B4").Value = 4755
$ws.Range("C4").Value = 'Ovaliderad'
$ws.Range("D4").Value = 'LC'
$ws.Range("E4").Value = 100857
$ws.Range("F4").Value = 'Robust tickgnagare'
$ws.Range("G4").Value = 'Dorcatoma robusta'
$ws.Range("H4").Value = 'Strand, 1938'
$ws.Range("I4").Value = "'1"
$ws.Range("J4").Value = 'ex.'
$ws.Range("K4").Value = 'imago/adult'
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = 'fönsterfälla'
$ws.Range("P4").Value = 'Åsums fure, delomr 19, 580 m NO om mc-banans ledningstorn, Sk'
$ws.Range("Q4").Value = 445824
$ws.Range("R4").Value = 6205171
$ws.Range("S4").Value = 10
$ws.Range("T4").Value = 'Skåne'
$ws.Range("U4").Value = 'Kristianstad'
$ws.Range("V4").Value = 'Skåne'
$ws.Range("W4").Value = 'Kristianstad'
$ws.Range("Y4").Value = "'2013-05-16"
$ws.Range("AA4").Value = "'2013-05-24"
$ws.Range("AD4").Value = $False
$ws.Range("AE4").Value = $False
$ws.Range("AG4").Value = $False
$ws.Range("AI4").Value = 'i gles tallskog'
$ws.Range("AO4").Value = 'på nydöd tall, delvis barklös'
$ws.Range("AQ4").Value = 'Nils Otto Nilsson'
$ws.Range("AR4").Value = 'NON 04741'
$ws.Range("AT4").ClearContents()
$ws.Range("AW4").Value = 'Nils Otto Nilsson'
$ws.Range("AX4").Value = 'Nils Otto Nilsson'
$ws.Range("AY4").Value = 'Åsums fure 2013'

# ---- Row 5 ----
$ws.Range("A5").Value = 112156964
$ws.Range("B5").Value = 39455
$ws.Range("C5").Value = 'Ovaliderad'
$ws.Range("D5").Value = 'NT'
$ws.Range("E5").Value = 102471
$ws.Range("F5").Value = 'Åkerväddsantennmal'
$ws.Range("G5").Value = 'Nemophora metallica'
$ws.Range("H5").Value = '(Poda, 1761)'
$ws.Range("I5").ClearContents()
$ws.Range("J5").ClearContents()
$ws.Range("K5").ClearContents()
$ws.Range("M5").Value = 'födosökande'
$ws.Range("N5").ClearContents()
$ws.Range("P5").Value = 'Åsums fure, delomr 19, 580 m NO om mc-banans ledningstorn, Sk'
$ws.Range("Q5").Value = 445828
$ws.Range("R5").Value = 6205165
$ws.Range("S5").Value = 10
$ws.Range("T5").Value = 'Skåne'
$ws.Range("U5").Value = 'Kristianstad'
$ws.Range("V5").Value = 'Skåne'
$ws.Range("W5").Value = 'Kristianstad'
$ws.Range("Y5").Value = "'2013-07-11"
$ws.Range("AA5").Value = "'2013-07-11"
$ws.Range("AD5").Value = $False
$ws.Range("AE5").Value = $False
$ws.Range("AG5").Value = $False
$ws.Range("AI5").Value = 'i gles tallskog'
$ws.Range("AO5").Value = 'på blmr av åkervädd'
$ws.Range("AQ5").ClearContents()
$ws.Range("AR5").ClearContents()
$ws.Range("AT5").ClearContents()
$ws.Range("AW5").Value = 'Nils Otto Nilsson'
$ws.Range("AX5").Value = 'Nils Otto Nilsson'
$ws.Range("AY5").Value = 'Åsums fure 2013'
